$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("backlog")

# --- Header row ---
$ws.Range("C1").Value = "feature"

# --- Row 2 ---
$ws.Range("C2").Value = "output"
$ws.Range("D2").Value = "output to excel"

# --- Row 3 ---
$ws.Range("C3").Value = "output"
$ws.Range("D3").Value = "output to screen"

# --- Row 4 ---
$ws.Range("C4").Value = "get CRS"
$ws.Range("D4").Value = "get_crs_type"
$ws.Range("F4").Value = "return the CRS or NONE"

# --- Row 5 ---
$ws.Range("C5").Value = "get CRS"
$ws.Range("D5").Value = "get_crs_wkid"
$ws.Range("F5").Value = "returns the coordinate reference system factory code (EPSG/WKID)"

# --- Row 6 ---
$ws.Range("C6").Value = "get CRS"
$ws.Range("D6").Value = "get_crs_name"
$ws.Range("E6").Value = "done"
$ws.Range("F6").Value = "return the CRS name"

# --- Row 7 ---
$ws.Range("C7").Value = "summary"
$ws.Range("D7").Value = "get_row_count"

# --- Row 8 ---
$ws.Range("C8").Value = "summary"
$ws.Range("D8").Value = "get_field_count"

# --- Row 9 ---
$ws.Range("C9").Value = "summary"
$ws.Range("D9").Value = "get_null_count(field)"
$ws.Range("F9").Value = "returns the count of nulls for a field"

# --- Update selected cell from A11 to E7 ---
$ws.Range("E7").Select()
